$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 18.31647966666667
$ws.Cells.Item(2, 8).Value = 54.949439
$ws.Cells.Item(2, 9).Value = 0.005487334033884006
$ws.Cells.Item(2, 10).Value = 0.005487334033884005
$ws.Cells.Item(2, 13).Value = 209.26237
$ws.Cells.Item(2, 14).Value = 627.78711
$ws.Cells.Item(2, 15).Value = 0.8127157202241573
$ws.Cells.Item(2, 16).Value = 0.8127157202241573
$ws.Cells.Item(2, 17).Value = 3832.949945103477
$ws.Cells.Item(2, 18).Value = 34496.54950593129
$ws.Cells.Item(2, 19).Value = 0.00445964263145857
$ws.Cells.Item(2, 20).Value = 0.004459642631458569
$ws.Cells.Item(3, 7).Value = 18.31647966666667
$ws.Cells.Item(3, 8).Value = 54.949439
$ws.Cells.Item(3, 9).Value = 0.005487334033884006
$ws.Cells.Item(3, 10).Value = 0.005487334033884005
$ws.Cells.Item(3, 13).Value = 0.9848756666666668
$ws.Cells.Item(3, 14).Value = 2.954627
$ws.Cells.Item(3, 15).Value = 0.003824977881910862
$ws.Cells.Item(3, 16).Value = 0.003824977881910862
$ws.Cells.Item(3, 17).Value = 18.03945512269478
$ws.Cells.Item(3, 18).Value = 162.355096104253
$ws.Cells.Item(3, 19).Value = 0.00002098893131026303
$ws.Cells.Item(3, 20).Value = 0.00002098893131026303
$ws.Cells.Item(4, 7).Value = 18.31647966666667
$ws.Cells.Item(4, 8).Value = 54.949439
$ws.Cells.Item(4, 9).Value = 0.005487334033884006
$ws.Cells.Item(4, 10).Value = 0.005487334033884005
$ws.Cells.Item(4, 13).Value = 1.763846666666667
$ws.Cells.Item(4, 14).Value = 5.291539999999999
$ws.Cells.Item(4, 15).Value = 0.006850280411451801
$ws.Cells.Item(4, 16).Value = 0.006850280411451801
$ws.Cells.Item(4, 17).Value = 32.30746160511777
$ws.Cells.Item(4, 18).Value = 290.76715444606
$ws.Cells.Item(4, 19).Value = 0.0000375897768434084
$ws.Cells.Item(4, 20).Value = 0.00003758977684340839
$ws.Cells.Item(5, 7).Value = 18.31647966666667
$ws.Cells.Item(5, 8).Value = 54.949439
$ws.Cells.Item(5, 9).Value = 0.005487334033884006
$ws.Cells.Item(5, 10).Value = 0.005487334033884005
$ws.Cells.Item(5, 13).Value = 45.474231
$ws.Cells.Item(5, 14).Value = 136.422693
$ws.Cells.Item(5, 15).Value = 0.1766090214824801
$ws.Cells.Item(5, 16).Value = 0.1766090214824801
$ws.Cells.Item(5, 17).Value = 832.9278274688029
$ws.Cells.Item(5, 18).Value = 7496.350447219225
$ws.Cells.Item(5, 19).Value = 0.0009691126942717645
$ws.Cells.Item(5, 20).Value = 0.0009691126942717644
$ws.Cells.Item(6, 9).Value = 0.9472399998689139
$ws.Cells.Item(6, 10).Value = 0.9472399998689137
$ws.Cells.Item(6, 13).Value = 209.26237
$ws.Cells.Item(6, 14).Value = 627.78711
$ws.Cells.Item(6, 15).Value = 0.8127157202241573
$ws.Cells.Item(6, 16).Value = 0.8127157202241573
$ws.Cells.Item(6, 17).Value = 661655.2743240779
$ws.Cells.Item(6, 18).Value = 5954897.4689167
$ws.Cells.Item(6, 19).Value = 0.769836838718595
$ws.Cells.Item(6, 20).Value = 0.7698368387185949
$ws.Cells.Item(7, 9).Value = 0.9472399998689139
$ws.Cells.Item(7, 10).Value = 0.9472399998689137
$ws.Cells.Item(7, 13).Value = 0.9848756666666668
$ws.Cells.Item(7, 14).Value = 2.954627
$ws.Cells.Item(7, 15).Value = 0.003824977881910862
$ws.Cells.Item(7, 16).Value = 0.003824977881910862
$ws.Cells.Item(7, 17).Value = 3114.024654329598
$ws.Cells.Item(7, 18).Value = 28026.22188896638
$ws.Cells.Item(7, 19).Value = 0.003623172048359844
$ws.Cells.Item(7, 20).Value = 0.003623172048359843
$ws.Cells.Item(8, 9).Value = 0.9472399998689139
$ws.Cells.Item(8, 10).Value = 0.9472399998689137
$ws.Cells.Item(8, 13).Value = 1.763846666666667
$ws.Cells.Item(8, 14).Value = 5.291539999999999
$ws.Cells.Item(8, 15).Value = 0.006850280411451801
$ws.Cells.Item(8, 16).Value = 0.006850280411451801
$ws.Cells.Item(8, 17).Value = 5577.010573372287
$ws.Cells.Item(8, 18).Value = 50193.09516035057
$ws.Cells.Item(8, 19).Value = 0.006488859616045627
$ws.Cells.Item(8, 20).Value = 0.006488859616045626
$ws.Cells.Item(9, 9).Value = 0.9472399998689139
$ws.Cells.Item(9, 10).Value = 0.9472399998689137
$ws.Cells.Item(9, 13).Value = 45.474231
$ws.Cells.Item(9, 14).Value = 136.422693
$ws.Cells.Item(9, 15).Value = 0.1766090214824801
$ws.Cells.Item(9, 16).Value = 0.1766090214824801
$ws.Cells.Item(9, 17).Value = 143782.490788867
$ws.Cells.Item(9, 18).Value = 1294042.417099803
$ws.Cells.Item(9, 19).Value = 0.1672911294859135
$ws.Cells.Item(9, 20).Value = 0.1672911294859134
$ws.Cells.Item(10, 7).Value = 155.6514383333333
$ws.Cells.Item(10, 8).Value = 466.954315
$ws.Cells.Item(10, 9).Value = 0.04663076369111781
$ws.Cells.Item(10, 10).Value = 0.0466307636911178
$ws.Cells.Item(10, 13).Value = 209.26237
$ws.Cells.Item(10, 14).Value = 627.78711
$ws.Cells.Item(10, 15).Value = 0.8127157202241573
$ws.Cells.Item(10, 16).Value = 0.8127157202241573
$ws.Cells.Item(10, 17).Value = 32571.98887954219
$ws.Cells.Item(10, 18).Value = 293147.8999158797
$ws.Cells.Item(10, 19).Value = 0.03789755469782929
$ws.Cells.Item(10, 20).Value = 0.03789755469782929
$ws.Cells.Item(11, 7).Value = 155.6514383333333
$ws.Cells.Item(11, 8).Value = 466.954315
$ws.Cells.Item(11, 9).Value = 0.04663076369111781
$ws.Cells.Item(11, 10).Value = 0.0466307636911178
$ws.Cells.Item(11, 13).Value = 0.9848756666666668
$ws.Cells.Item(11, 14).Value = 2.954627
$ws.Cells.Item(11, 15).Value = 0.003824977881910862
$ws.Cells.Item(11, 16).Value = 0.003824977881910862
$ws.Cells.Item(11, 17).Value = 153.2973140961673
$ws.Cells.Item(11, 18).Value = 1379.675826865505
$ws.Cells.Item(11, 19).Value = 0.0001783616397351378
$ws.Cells.Item(11, 20).Value = 0.0001783616397351377
$ws.Cells.Item(12, 7).Value = 155.6514383333333
$ws.Cells.Item(12, 8).Value = 466.954315
$ws.Cells.Item(12, 9).Value = 0.04663076369111781
$ws.Cells.Item(12, 10).Value = 0.0466307636911178
$ws.Cells.Item(12, 13).Value = 1.763846666666667
$ws.Cells.Item(12, 14).Value = 5.291539999999999
$ws.Cells.Item(12, 15).Value = 0.006850280411451801
$ws.Cells.Item(12, 16).Value = 0.006850280411451801
$ws.Cells.Item(12, 17).Value = 274.5452706661222
$ws.Cells.Item(12, 18).Value = 2470.9074359951
$ws.Cells.Item(12, 19).Value = 0.0003194338070843022
$ws.Cells.Item(12, 20).Value = 0.0003194338070843022
$ws.Cells.Item(13, 7).Value = 155.6514383333333
$ws.Cells.Item(13, 8).Value = 466.954315
$ws.Cells.Item(13, 9).Value = 0.04663076369111781
$ws.Cells.Item(13, 10).Value = 0.0466307636911178
$ws.Cells.Item(13, 13).Value = 45.474231
$ws.Cells.Item(13, 14).Value = 136.422693
$ws.Cells.Item(13, 15).Value = 0.1766090214824801
$ws.Cells.Item(13, 16).Value = 0.1766090214824801
$ws.Cells.Item(13, 17).Value = 7078.129462252255
$ws.Cells.Item(13, 18).Value = 63703.16516027028
$ws.Cells.Item(13, 19).Value = 0.008235413546469079
$ws.Cells.Item(13, 20).Value = 0.008235413546469075
$ws.Cells.Item(14, 7).Value = 2.142642
$ws.Cells.Item(14, 8).Value = 6.427926
$ws.Cells.Item(14, 9).Value = 0.0006419024060843985
$ws.Cells.Item(14, 10).Value = 0.0006419024060843984
$ws.Cells.Item(14, 13).Value = 209.26237
$ws.Cells.Item(14, 14).Value = 627.78711
$ws.Cells.Item(14, 15).Value = 0.8127157202241573
$ws.Cells.Item(14, 16).Value = 0.8127157202241573
$ws.Cells.Item(14, 17).Value = 448.37434298154
$ws.Cells.Item(14, 18).Value = 4035.36908683386
$ws.Cells.Item(14, 19).Value = 0.0005216841762745014
$ws.Cells.Item(14, 20).Value = 0.0005216841762745013
$ws.Cells.Item(15, 7).Value = 2.142642
$ws.Cells.Item(15, 8).Value = 6.427926
$ws.Cells.Item(15, 9).Value = 0.0006419024060843985
$ws.Cells.Item(15, 10).Value = 0.0006419024060843984
$ws.Cells.Item(15, 13).Value = 0.9848756666666668
$ws.Cells.Item(15, 14).Value = 2.954627
$ws.Cells.Item(15, 15).Value = 0.003824977881910862
$ws.Cells.Item(15, 16).Value = 0.003824977881910862
$ws.Cells.Item(15, 17).Value = 2.110235968178
$ws.Cells.Item(15, 18).Value = 18.992123713602
$ws.Cells.Item(15, 19).Value = 0.000002455262505618189
$ws.Cells.Item(15, 20).Value = 0.000002455262505618188
$ws.Cells.Item(16, 7).Value = 2.142642
$ws.Cells.Item(16, 8).Value = 6.427926
$ws.Cells.Item(16, 9).Value = 0.0006419024060843985
$ws.Cells.Item(16, 10).Value = 0.0006419024060843984
$ws.Cells.Item(16, 13).Value = 1.763846666666667
$ws.Cells.Item(16, 14).Value = 5.291539999999999
$ws.Cells.Item(16, 15).Value = 0.006850280411451801
$ws.Cells.Item(16, 16).Value = 0.006850280411451801
$ws.Cells.Item(16, 17).Value = 3.77929194956
$ws.Cells.Item(16, 18).Value = 34.01362754604
$ws.Cells.Item(16, 19).Value = 0.000004397211478463734
$ws.Cells.Item(16, 20).Value = 0.000004397211478463734
$ws.Cells.Item(17, 7).Value = 2.142642
$ws.Cells.Item(17, 8).Value = 6.427926
$ws.Cells.Item(17, 9).Value = 0.0006419024060843985
$ws.Cells.Item(17, 10).Value = 0.0006419024060843984
$ws.Cells.Item(17, 13).Value = 45.474231
$ws.Cells.Item(17, 14).Value = 136.422693
$ws.Cells.Item(17, 15).Value = 0.1766090214824801
$ws.Cells.Item(17, 16).Value = 0.1766090214824801
$ws.Cells.Item(17, 17).Value = 97.43499725830199
$ws.Cells.Item(17, 18).Value = 876.9149753247179
$ws.Cells.Item(17, 19).Value = 0.0001133657558258152
$ws.Cells.Item(17, 20).Value = 0.0001133657558258152
